$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 509.53333
$ws.Range("I19").Value = 343.5
$ws.Range("J19").Value = 699.2857
$ws.Range("K19").Value = 343.5
$ws.Range("L19").Value = 699.2857
$ws.Range("M19").Value = -168.5
$ws.Range("N19").Value = -1049.2857

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1526.5555
$ws.Range("I53").Value = 647.8
$ws.Range("J53").Value = 2625
$ws.Range("K53").Value = 647.8
$ws.Range("L53").Value = 2625
$ws.Range("M53").Value = -10.79999999999995
$ws.Range("N53").Value = -3899

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2395.3333
$ws.Range("I86").Value = 1091.75
$ws.Range("K86").Value = 1091.75
$ws.Range("M86").Value = 31.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 2395.3333
$ws.Range("I89").Value = 1091.75
$ws.Range("K89").Value = 5458.75
$ws.Range("M89").Value = 157.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 2562.5
$ws.Range("I101").Value = 2000
$ws.Range("J101").Value = 3125
$ws.Range("K101").Value = 6000
$ws.Range("L101").Value = 9375
$ws.Range("M101").Value = -4378
$ws.Range("N101").Value = -12619

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1529.7
$ws.Range("I129").Value = 637
$ws.Range("J129").Value = 1752.875
$ws.Range("K129").Value = 1911
$ws.Range("L129").Value = 5258.625
$ws.Range("M129").Value = 3089
$ws.Range("N129").Value = -15258.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2543.9375
$ws.Range("I137").Value = 1500
$ws.Range("J137").Value = 2613.5334
$ws.Range("K137").Value = 4500
$ws.Range("L137").Value = 7840.600199999999
$ws.Range("M137").Value = -1950
$ws.Range("N137").Value = -12940.6002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1254077.4
$ws.Range("I138").Value = 2405.5
$ws.Range("J138").Value = 2505749.2
$ws.Range("K138").Value = 7216.5
$ws.Range("L138").Value = 7517247.600000001
$ws.Range("M138").Value = -2076.5
$ws.Range("N138").Value = -7527527.600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 8239.6
$ws.Range("J74").Value = 5300
$ws.Range("L74").Value = 5300
$ws.Range("N74").Value = -7048

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 8239.6
$ws.Range("J77").Value = 5300
$ws.Range("L77").Value = 26500
$ws.Range("N77").Value = -35236

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1716212.2
$ws.Range("I86").Value = 2477043
$ws.Range("K86").Value = 2477043
$ws.Range("M86").Value = -2475920

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1716212.2
$ws.Range("I89").Value = 2477043
$ws.Range("K89").Value = 12385215
$ws.Range("M89").Value = -12379599

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3717.625
$ws.Range("I105").Value = 2606.8462
$ws.Range("J105").Value = 4477.6313
$ws.Range("K105").Value = 2606.8462
$ws.Range("L105").Value = 4477.6313
$ws.Range("M105").Value = -859.8462
$ws.Range("N105").Value = -7971.6313

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5497.2856
$ws.Range("I31").Value = 4261
$ws.Range("J31").Value = 5943.722
$ws.Range("K31").Value = 4261
$ws.Range("L31").Value = 5943.722
$ws.Range("M31").Value = -3966
$ws.Range("N31").Value = -6533.722

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5497.2856
$ws.Range("I34").Value = 4261
$ws.Range("J34").Value = 5943.722
$ws.Range("K34").Value = 4261
$ws.Range("L34").Value = 5943.722
$ws.Range("M34").Value = -4059
$ws.Range("N34").Value = -6347.722

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2717.5
$ws.Range("I58").Value = 3435
$ws.Range("K58").Value = 3435
$ws.Range("M58").Value = -3232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2717.5
$ws.Range("I136").Value = 3435
$ws.Range("K136").Value = 10305
$ws.Range("M136").Value = -7755

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 799255.9399999999
$ws.Range("I56").Value = 799255.9399999999
$ws.Range("K56").Value = 799255.9399999999
$ws.Range("M56").Value = -798725.9399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 4155
$ws.Range("I59").Value = 850
$ws.Range("J59").Value = 7460
$ws.Range("K59").Value = 2550
$ws.Range("L59").Value = 22380
$ws.Range("M59").Value = -2010
$ws.Range("N59").Value = -23460

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 7183.1665
$ws.Range("J68").Value = 9899.75
$ws.Range("L68").Value = 29699.25
$ws.Range("N68").Value = -31321.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 7183.1665
$ws.Range("J71").Value = 9899.75
$ws.Range("L71").Value = 89097.75
$ws.Range("N71").Value = -97209.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1526.625
$ws.Range("J107").Value = 2259.8
$ws.Range("L107").Value = 6779.400000000001
$ws.Range("N107").Value = -10619.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 904.25
$ws.Range("I122").Value = 558.25
$ws.Range("J122").Value = 1077.25
$ws.Range("K122").Value = 5024.25
$ws.Range("L122").Value = 9695.25
$ws.Range("M122").Value = -2574.25
$ws.Range("N122").Value = -14595.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 38895780
$ws.Range("I137").Value = 50004304
$ws.Range("J137").Value = 33341520
$ws.Range("K137").Value = 150012912
$ws.Range("L137").Value = 100024560
$ws.Range("M137").Value = -150007812
$ws.Range("N137").Value = -100034760

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8063.864
$ws.Range("I70").Value = 5757.727
$ws.Range("J70").Value = 10370
$ws.Range("K70").Value = 5757.727
$ws.Range("L70").Value = 10370
$ws.Range("M70").Value = -5487.727
$ws.Range("N70").Value = -10910

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 8063.864
$ws.Range("I73").Value = 5757.727
$ws.Range("J73").Value = 10370
$ws.Range("K73").Value = 5757.727
$ws.Range("L73").Value = 10370
$ws.Range("M73").Value = -4821.727
$ws.Range("N73").Value = -12242

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 37079120
$ws.Range("I122").Value = 52689480
$ws.Range("J122").Value = 4515.625
$ws.Range("K122").Value = 158068440
$ws.Range("L122").Value = 13546.875
$ws.Range("M122").Value = -158065990
$ws.Range("N122").Value = -18446.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 48230.832
$ws.Range("J123").Value = 48230.832
$ws.Range("L123").Value = 48230.832
$ws.Range("N123").Value = -53130.832

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 141350.25
$ws.Range("J139").Value = 141350.25
$ws.Range("L139").Value = 141350.25
$ws.Range("N139").Value = -151630.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4831.365
$ws.Range("J132").Value = 5222.1875
$ws.Range("L132").Value = 15666.5625
$ws.Range("N132").Value = -20726.5625
